$wb = $excel.ActiveWorkbook

# Sheet "Folder Inventory" - C2: Last Updated timestamp
$wsInventory = $wb.Worksheets.Item("Folder Inventory")
$wsInventory.Range("C2").Value = "2025-06-12 17:36:45 +0530"

# Sheet "Metadata" - B3: Generated On, B5: Workflow Run
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B3").Value = "2025-06-12 12:07:04 UTC"
# Workflow Run must stay a text value ("14"), not be auto-converted to a number.
$wsMetadata.Range("B5").NumberFormat = "@"
$wsMetadata.Range("B5").Value = "14"
$wsMetadata.Range("B5").ClearFormats()

# Sheet "Summary" - B5: Most Recent Update
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-12 17:36:45 +0530"
